# "ticks per meter calculations" - update Pulses/Rev-to-Ticks/Meter inputs
# and add left/right wheel-circumference rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Inputs section (rows 11-14) -----------------------------------------
# Diameter (inches) input E12: 5 -> 4.986
$ws.Range("E12").Value = 4.986

# Pulses/Rev B13: 512 -> 509.5
$ws.Range("B13").Value = 509.5

# --- New rows 16 & 17: per-side readings ----------------------------------
# A16/A17 pick up the same highlight style used by the other Inputs labels
# (A12:A14), so copy that formatting across before writing the new labels.
$ws.Range("A14").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = "right"
$ws.Range("B16").Value = 504.75

$ws.Range("A17").Value = "left"
$ws.Range("B17").Value = 508.5

# --- Selection -------------------------------------------------------------
$ws.Range("B14").Select()
